# Remove the " (Remote)" qualifier that follows "Software Team Lead" and
# relocate the document's "_GoBack" bookmark (Word's "last edit location"
# marker) from its old spot after ", IIT Bombay" to the point where the
# text was deleted - exactly what Word itself does when you delete text.

$d = $word.ActiveDocument

# --- 1. Drop the stray "_GoBack" bookmark currently sitting after
#        ", IIT Bombay" (it will be re-created at the real edit point).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Locate "Software Team Lead" - the new bookmark belongs right
#        after it, which is also where the " (Remote)" run begins.
$lead = $d.Content
$lead.Find.Execute("Software Team Lead", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $lead.End

# --- 3. Locate the " (Remote)" run that needs to disappear.
$remote = $d.Content
$remote.Find.Execute(" (Remote)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# --- 4. Bracket the text to be removed with bookmarks *before* deleting
#        it. Both neighbouring runs share identical run formatting, so a
#        plain delete would let Word coalesce them into the bookmark's
#        neighbouring runs; the bookmarks keep them distinct, matching
#        the original run layout exactly.
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))
$guardName = "ZZGuard_TempSplit"
$d.Bookmarks.Add($guardName, $d.Range($remote.End, $remote.End))

# --- 5. Delete the " (Remote)" text itself.
$d.Range($splitPoint, $remote.End).Delete()

# --- 6. Remove the temporary guard bookmark, leaving only "_GoBack".
$d.Bookmarks($guardName).Delete()
